# Auto commit at 2026-01-06  9:09:13.91
# Append the new daily readings (2026-01-05, serial 46027) for both
# charging stations to the bottom of the data table on Sheet1, and
# update the view so the newly added rows are visible/selected the
# way the author left the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New row 10: 四方坪站 (station 1) for 2026-01-05 -----------------
$ws.Cells.Item(10, 1).Value = 46027
$ws.Cells.Item(10, 2).Value = "四方坪站"
$ws.Cells.Item(10, 3).Value = 13599.47
$ws.Cells.Item(10, 4).Value = 9534.6
$ws.Cells.Item(10, 5).Value = 3014.77
$ws.Cells.Item(10, 6).Value = 574

# --- New row 11: 高岭站 (station 2) for 2026-01-05 --------------------
$ws.Cells.Item(11, 1).Value = 46027
$ws.Cells.Item(11, 2).Value = "高岭站"
$ws.Cells.Item(11, 3).Value = 5911.55
$ws.Cells.Item(11, 4).Value = 4914.84
$ws.Cells.Item(11, 5).Value = 1637.77
$ws.Cells.Item(11, 6).Value = 198

# --- Update the window/view so the new rows are scrolled into view --
$win = $excel.ActiveWindow
$win.ScrollRow = 10
$win.ScrollColumn = 1

# --- Leave the same cell selected as the author did ------------------
$ws.Range("H15").Select()
